# CIV-11205 updated templates with court full address
#
# The "This order is made by <<judgeNameTitle>> on <<submittedOn>> at
# <<courtName>>." paragraph needs its trailing "<<courtName>>" merge-field
# replaced with the full court-site address merge-fields (siteName,
# address, postcode).

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    ">> at <<courtName>>.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ">> at <<siteName>> - <<address>> - <<postcode>>.",
    2
)
